$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting existing rows 250:325 down to 251:326
$ws.Rows("250:250").Insert()

# Populate the newly inserted row 250 with its data
$ws.Range("A250").Value = 4
$ws.Range("B250").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C250").Value = "Los Lagos"
$ws.Range("D250").Value = 44988
$ws.Range("E250").Value = 10
$ws.Range("F250").Value = 100112039
$ws.Range("G250").Value = "Ciboulette"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 220
$ws.Range("K250").Value = 3500
$ws.Range("L250").Value = 3500
$ws.Range("M250").Value = 3500
$ws.Range("N250").Value = "$/docena de atados"
$ws.Range("O250").Value = "Región Metropolitana"
$ws.Range("P250").Value = 1167
$ws.Range("Q250").Value = 3
$ws.Range("R250").Value = "Hortaliza"
